$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.04332866666666666
$ws.Cells.Item(2, 8).Value = 0.129986
$ws.Cells.Item(2, 9).Value = 0.2246397599897691
$ws.Cells.Item(2, 10).Value = 0.2246397599897691
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.007258333333333333
$ws.Cells.Item(2, 14).Value = 0.021775
$ws.Cells.Item(2, 15).Value = 0.000328667160253549
$ws.Cells.Item(2, 16).Value = 0.000328667160253549
$ws.Cells.Item(2, 17).Value = 0.0003144939055555555
$ws.Cells.Item(2, 18).Value = 0.00283044515
$ws.Cells.Item(2, 19).Value = 0.00007383171199587624
$ws.Cells.Item(2, 20).Value = 0.00007383171199587623
$ws.Cells.Item(3, 7).Value = 0.04332866666666666
$ws.Cells.Item(3, 8).Value = 0.129986
$ws.Cells.Item(3, 9).Value = 0.2246397599897691
$ws.Cells.Item(3, 10).Value = 0.2246397599897691
$ws.Cells.Item(3, 15).Value = 0.7778551418094273
$ws.Cells.Item(3, 16).Value = 0.7778551418094272
$ws.Cells.Item(3, 17).Value = 0.7443113614253333
$ws.Cells.Item(3, 18).Value = 6.698802252828
$ws.Cells.Item(3, 19).Value = 0.1747371923628776
$ws.Cells.Item(3, 20).Value = 0.1747371923628776
$ws.Cells.Item(4, 7).Value = 0.04332866666666666
$ws.Cells.Item(4, 8).Value = 0.129986
$ws.Cells.Item(4, 9).Value = 0.2246397599897691
$ws.Cells.Item(4, 10).Value = 0.2246397599897691
$ws.Cells.Item(4, 13).Value = 4.898620999999999
$ws.Cells.Item(4, 14).Value = 14.695863
$ws.Cells.Item(4, 15).Value = 0.2218161910303192
$ws.Cells.Item(4, 16).Value = 0.2218161910303192
$ws.Cells.Item(4, 17).Value = 0.2122507164353333
$ws.Cells.Item(4, 18).Value = 1.910256447918
$ws.Cells.Item(4, 19).Value = 0.04982873591489569
$ws.Cells.Item(4, 20).Value = 0.04982873591489569
$ws.Cells.Item(5, 9).Value = 0.5955530362469368
$ws.Cells.Item(5, 10).Value = 0.5955530362469369
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.007258333333333333
$ws.Cells.Item(5, 14).Value = 0.021775
$ws.Cells.Item(5, 15).Value = 0.000328667160253549
$ws.Cells.Item(5, 16).Value = 0.000328667160253549
$ws.Cells.Item(5, 17).Value = 0.0008337695888888889
$ws.Cells.Item(5, 18).Value = 0.007503926300000001
$ws.Cells.Item(5, 19).Value = 0.0001957387252036597
$ws.Cells.Item(5, 20).Value = 0.0001957387252036597
$ws.Cells.Item(6, 9).Value = 0.5955530362469368
$ws.Cells.Item(6, 10).Value = 0.5955530362469369
$ws.Cells.Item(6, 15).Value = 0.7778551418094273
$ws.Cells.Item(6, 16).Value = 0.7778551418094272
$ws.Cells.Item(6, 19).Value = 0.463253991464896
$ws.Cells.Item(6, 20).Value = 0.463253991464896
$ws.Cells.Item(7, 9).Value = 0.5955530362469368
$ws.Cells.Item(7, 10).Value = 0.5955530362469369
$ws.Cells.Item(7, 13).Value = 4.898620999999999
$ws.Cells.Item(7, 14).Value = 14.695863
$ws.Cells.Item(7, 15).Value = 0.2218161910303192
$ws.Cells.Item(7, 16).Value = 0.2218161910303192
$ws.Cells.Item(7, 17).Value = 0.5627078600173333
$ws.Cells.Item(7, 18).Value = 5.064370740156
$ws.Cells.Item(7, 19).Value = 0.1321033060568372
$ws.Cells.Item(7, 20).Value = 0.1321033060568372
$ws.Cells.Item(8, 7).Value = 0.03468133333333333
$ws.Cells.Item(8, 8).Value = 0.104044
$ws.Cells.Item(8, 9).Value = 0.1798072037632941
$ws.Cells.Item(8, 10).Value = 0.1798072037632941
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.007258333333333333
$ws.Cells.Item(8, 14).Value = 0.021775
$ws.Cells.Item(8, 15).Value = 0.000328667160253549
$ws.Cells.Item(8, 16).Value = 0.000328667160253549
$ws.Cells.Item(8, 17).Value = 0.0002517286777777777
$ws.Cells.Item(8, 18).Value = 0.0022655581
$ws.Cells.Item(8, 19).Value = 0.00005909672305401311
$ws.Cells.Item(8, 20).Value = 0.0000590967230540131
$ws.Cells.Item(9, 7).Value = 0.03468133333333333
$ws.Cells.Item(9, 8).Value = 0.104044
$ws.Cells.Item(9, 9).Value = 0.1798072037632941
$ws.Cells.Item(9, 10).Value = 0.1798072037632941
$ws.Cells.Item(9, 15).Value = 0.7778551418094273
$ws.Cells.Item(9, 16).Value = 0.7778551418094272
$ws.Cells.Item(9, 17).Value = 0.5957651692346667
$ws.Cells.Item(9, 18).Value = 5.361886523112
$ws.Cells.Item(9, 19).Value = 0.1398639579816537
$ws.Cells.Item(9, 20).Value = 0.1398639579816537
$ws.Cells.Item(10, 7).Value = 0.03468133333333333
$ws.Cells.Item(10, 8).Value = 0.104044
$ws.Cells.Item(10, 9).Value = 0.1798072037632941
$ws.Cells.Item(10, 10).Value = 0.1798072037632941
$ws.Cells.Item(10, 13).Value = 4.898620999999999
$ws.Cells.Item(10, 14).Value = 14.695863
$ws.Cells.Item(10, 15).Value = 0.2218161910303192
$ws.Cells.Item(10, 16).Value = 0.2218161910303192
$ws.Cells.Item(10, 17).Value = 0.1698907077746667
$ws.Cells.Item(10, 18).Value = 1.529016369972
$ws.Cells.Item(10, 19).Value = 0.03988414905858637
$ws.Cells.Item(10, 20).Value = 0.03988414905858637
